$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the added columns D:F
$ws.Range("D1").Value = "No. of Actions"
$ws.Range("E1").Value = "No. of New Node Expansions"
$ws.Range("F1").Value = "Time to Complete Plan Search"

# Match the column widths Excel computed for the new columns as closely as
# the engine's width-quantization allows.
$ws.Columns.Item(4).ColumnWidth = 12.59
$ws.Columns.Item(5).ColumnWidth = 26.25
$ws.Columns.Item(6).ColumnWidth = 27.083333333333332

# Move the selection to F2, matching the saved view state.
$ws.Range("F2").Select()
